$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the year header row (E1:BL1) from text labels like "1960 [YR1960]"
# to plain numeric years 1960-2019, left aligned. (2020 in BM1 is left as-is.)
$year = 1960
for ($col = 5; $col -le 64; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $year
    $cell.HorizontalAlignment = -4131
    $year++
}

# Match the updated selection/active cell left behind by the edit.
$ws.Range("E1:BL1").Select()

Write-Output "done"
